# Edit described by the diff:
#  1. Worksheet "ODI Batting" (sheet2): clear the empty INNING_NUMBER
#     placeholder cells in column B for rows 2,3,4,6,8,9 (rows that have
#     no recorded inning number), leaving them as true blanks.
#  2. Add a brand-new worksheet "ODI Batting Extra" (sheet4) at the end
#     of the workbook with MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 /
#     PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH columns.

$wb = $excel.ActiveWorkbook

# --- 1. Clean up "ODI Batting" ---------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("B2").ClearContents()
$batting.Range("B3").ClearContents()
$batting.Range("B4").ClearContents()
$batting.Range("B6").ClearContents()
$batting.Range("B8").ClearContents()
$batting.Range("B9").ClearContents()

# --- 2. Add "ODI Batting Extra" worksheet at the end ------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$extra.Name = "ODI Batting Extra"

# Header row (bold, centered, bordered - same look as the other sheets)
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $extra.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$headerRange = $extra.Range("A1:F1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data rows: MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$data = @(
    @("4408", $null, $null, $null, $null, "NO"),
    @("4426", $null, $null, $null, $null, "NO"),
    @("4427", 11,    "1",   "0",   "3.66%", "NO"),
    @("4428", 10,    $null, $null, $null, "YES"),
    @("4472", 10,    "1",   "0",   "3.24%", "NO"),
    @("4473", 10,    $null, $null, $null, "NO"),
    @("4476", 10,    $null, $null, $null, "NO"),
    @("4713", $null, $null, $null, $null, $null)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $r + 2
    $rec = $data[$r]

    # MATCH_CODE - always text
    $a = $extra.Cells.Item($row, 1)
    $a.NumberFormat = "@"
    $a.Value = $rec[0]

    # BATTING_POSITION - numeric when present
    if ($null -ne $rec[1]) {
        $extra.Cells.Item($row, 2).Value = $rec[1]
    }

    # NUM_4 - text when present
    if ($null -ne $rec[2]) {
        $c = $extra.Cells.Item($row, 3)
        $c.NumberFormat = "@"
        $c.Value = $rec[2]
    }

    # NUM_6 - text when present
    if ($null -ne $rec[3]) {
        $d = $extra.Cells.Item($row, 4)
        $d.NumberFormat = "@"
        $d.Value = $rec[3]
    }

    # PERCENT_RUNS_OF_TOTAL - text when present
    if ($null -ne $rec[4]) {
        $e = $extra.Cells.Item($row, 5)
        $e.NumberFormat = "@"
        $e.Value = $rec[4]
    }

    # MAN_OF_MATCH - text when present
    if ($null -ne $rec[5]) {
        $f = $extra.Cells.Item($row, 6)
        $f.NumberFormat = "@"
        $f.Value = $rec[5]
    }
}
